# Commit message: "renamed everything from 'log type' to 'reporting type'"
#
# In this workbook that shows up in two places:
#   1. The worksheet tab "Log_Type_Statistics".
#   2. The column header text "log_type" in cell A1 of that same sheet.
#
# (The second line of the commit message, about the sampling-delay symbol,
# refers to other files in the same commit and has no corresponding content
# in this spreadsheet.)

$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("Log_Type_Statistics")

# Update the header text before renaming the sheet, then rename the tab.
$ws.Range("A1").Value = "reporting_type"
$ws.Name = "Reporting_Type_Statistics"
